$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt 28; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
}

$ws.Range("A5:A32").Select() | Out-Null
